$d = $word.ActiveDocument

# --- Change 1: add ru-RU language tagging to the "Цель курсовой работы" heading paragraph ---
$headingXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="1"/>
    <w:rPr>
      <w:color w:val="auto"/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:color w:val="auto"/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:lastRenderedPageBreak/>
    <w:t>Цель курсовой работы</w:t>
  </w:r>
</w:p>'

$found = $false
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("Цель курсовой работы")) {
        $p.Range.InsertXML($headingXml)
        $found = $true
        break
    }
}
if (-not $found) {
    throw "heading paragraph not found"
}

# --- Change 2: rewrite the two list-item paragraphs about the stop-location data store ---
$paraAXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="aa"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>'
$paraBXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="aa"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="1"/>
    </w:numPr>
    <w:rPr>
      <w:b/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
  </w:pPr>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve">Данные </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>об остановках</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:b/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> общественного транспорта. </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve">Для каждой остановки требуется хранить ее </w:t>
  </w:r>
  <w:r>
    <w:t>id</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>географическое местоположение (в виде координат)</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> и название</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> Для обеспечения быстродействия запросов </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>типа:</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>найти остановки в некотором районе, или найти ближайшую остановку,</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>или найти остановку по имени и т.д.</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> будем использовать колоночную</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> СУБД</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t xml:space="preserve"> (</w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:u w:val="single"/>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>Сassandra</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:lang w:val="ru-RU"/>
    </w:rPr>
    <w:t>)</w:t>
  </w:r>
  <w:r>
    <w:t>.</w:t>
  </w:r>
</w:p>'

$pA = $null
$pB = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.IndexOf("Redis") -ge 0) {
        $pA = $p
    }
}
if ($null -eq $pA) {
    throw "Redis paragraph not found"
}
$pB = $pA.Next()
if ($null -eq $pB) {
    throw "paragraph after Redis paragraph not found"
}

# Rewrite paragraph B (Расписание ...) first so paragraph A's range is unaffected.
$pB.Range.InsertXML($paraBXml)
$pA.Range.InsertXML($paraAXml)

# --- Change 3: drop the old _GoBack bookmark that used to sit after the diagram image ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

Write-Output "done"
